# Applies the commit "organized repo, added report" to programAnalysis.xlsx
#
# The edit adds a 5th measurement row (row 18) for the <64, 10000, 20> block
# on "Foglio1", which in turn changes the dependent average/summary formulas
# in row 21 (and the formulas in I16/J16 that depend on row 21) through
# normal Excel recalculation. It also updates the sheet's saved cell
# selection to J16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Foglio1")

# New data row for the <64, 10000, 20> experiment set (row 18, filling the
# previously-empty slot between row 17 and the summary row 21).
$ws.Range("B18").Value = 1869.93
$ws.Range("C18").Value = 1992000
$ws.Range("D18").Value = 1974274
$ws.Range("E18").Value = 2.39
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 0

# Row 21 (averages) and I16/J16 (which both read back from row 21) are
# driven by existing formulas, so they recompute automatically once row 18
# has values.

# Leave the saved selection on J16, matching the author's last position.
[void]$ws.Range("J16").Select()
